# Update cryptos list values to match the latest scrape (coinranking.com).
# Each assignment is prefixed with a leading apostrophe so that Excel
# stores numeric-looking / percent strings as literal text (matching the
# original inlineStr cells) instead of coercing them to numbers; the
# ClearFormats() call afterwards drops the quote-prefix marker Excel adds
# so the cell keeps its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '51.881.04'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'" + '  +0.26%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'" + '2.939.91'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'" + '  +4.01%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'" + '  +0.05%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'" + '352.51'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'" + '  +0.68%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'" + '112.15'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'" + '  -0.52%  '
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'" + '  +0.53%  '
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'" + '  +0.06%  '
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'" + '0.626'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'" + '  +1.27%  '
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'" + '39.39'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'" + '  -1.85%  '
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'" + '0.0891'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'" + '  +5.11%  '
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'" + '  +1.15%  '
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'" + '  +0.00%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'" + '7.86'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'" + '  +1.27%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'" + '3.400.75'
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'" + '2.936.50'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'" + '  +4.00%  '
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = "'" + '  +0.78%  '
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'" + '51.980.43'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'" + '  +0.39%  '
$ws.Range("E18").ClearFormats()
$ws.Range("B19").Value = "'" + 'ImmutableX'
$ws.Range("B19").ClearFormats()
$ws.Range("C19").Value = "'" + 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C19").ClearFormats()
$ws.Range("D19").Value = "'" + '3.33'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'" + '  -3.60%  '
$ws.Range("E19").ClearFormats()
$ws.Range("B20").Value = "'" + 'Uniswap'
$ws.Range("B20").ClearFormats()
$ws.Range("C20").Value = "'" + 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C20").ClearFormats()
$ws.Range("D20").Value = "'" + '7.65'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'" + '  +0.81%  '
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'" + '14.39'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'" + '  +7.43%  '
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'" + '0.0₃0988'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'" + '  +1.70%  '
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'" + '71.25'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'" + '  +1.22%  '
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'" + '270.00'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'" + '  +0.42%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'" + '2.79'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'" + '  +1.44%  '
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'" + '  +9.69%  '
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'" + '26.97'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'" + '  +2.78%  '
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'" + '  +0.19%  '
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'" + '7.44'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'" + '  +17.30%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'" + '0.109'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'" + '  +21.54%  '
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'" + '  +0.82%  '
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'" + '37.65'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'" + '  -2.42%  '
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "'" + '  +0.25%  '
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'" + '6.18'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'" + '  +9.98%  '
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'" + '  +0.20%  '
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'" + '  +1.24%  '
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'" + '  -0.20%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'" + '3.31'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'" + '  +3.30%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'" + '18.86'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'" + '  -0.31%  '
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'" + '  +1.93%  '
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'" + '2.70'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'" + '  +7.05%  '
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'" + '  +1.76%  '
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'" + '23.45'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'" + '  +6.63%  '
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'" + '  -1.00%  '
$ws.Range("E44").ClearFormats()
$ws.Range("E45").Value = "'" + '  +0.40%  '
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'" + '  +1.63%  '
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'" + '2.171.94'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'" + '  +0.30%  '
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'" + '112.43'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'" + '  -8.24%  '
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'" + '0.247'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'" + '  +0.25%  '
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'" + '0.0343'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'" + '  +10.78%  '
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'" + '  -1.22%  '
$ws.Range("E51").ClearFormats()
